$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 2025 metadata refresh for indicator 7.1.1: updated indicator wording,
# organization department, and contact details (uploaded new version of
# the metadata sheet).

# Indicator text now mentions "stable" access to electricity
$ws.Range("B4").Value = "7.1.1. Доля населения, имеющего стабильный доступ к электроэнергии "

# Organization: department renamed from "Отдел" to "Управление"
$ws.Range("B6").Value = "Национальный статистический комитет КР" + [char]10 + "(Управление статистики домашних хозяйств)"

# New contact person
$ws.Range("B7").Value = "Калымбетова Ы.И."

# New phone number
$ws.Range("B9").Value = "(0312) 32 46 55"

# New organization website
$ws.Range("B10").Value = "www.stat.gov.kg"

# New contact e-mail (edited last, matching the author's editing order)
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "

# Leave the active selection on the last-edited cell, as in the source file
$ws.Range("B8").Select()
